$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "The methodology ... real-world." paragraph gains a trailing space after
#    the final period (content is otherwise unchanged - the remaining diff
#    hunks for earlier paragraphs are pure run-merges that do not alter the
#    visible text, so no further edits are required for them).
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("complexities inherent in real-world.", $true, $true, $false, $false, $false, $true, 1, $false, "complexities inherent in real-world. ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Two new paragraphs are appended at the end of the document:
#      - "Notes" (bold heading)
#      - a body paragraph describing use of Google Collab / Gemini
#    Both empty paragraphs are created first (while still inheriting the
#    surrounding non-bold formatting), then text + bold formatting are
#    applied individually so the body paragraph is not left bold.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$lastPara.Range.InsertParagraphAfter()

$notesPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$bodyPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$notesPara.Range.Text = "Notes"
$bodyPara.Range.Text = "I used Google Collab for notebook-based coding and leveraged generative AI tool- Gemini to assist with code auto-completion and syntax suggestions, ensuring all implementations were written and customized by me."

$notesPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$notesPara.Range.Bold = 1
$notesPara.Range.BoldBi = 1
